$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.776.39"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "3.835.85"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'600.68"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'161.63"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").Value = "3.834.53"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'36.80"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "4.485.41"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "3.811.11"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "68.924.99"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'7.50"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'11.32"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "'484.07"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "'83.97"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D29").Value = "'9.95"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "'7.91"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").Value = "3.988.85"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").Value = "'32.04"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "3.786.56"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "'1.02"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "'0.139"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "'5.88"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D42").Value = "'436.71"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").Value = "'2.95"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "'48.48"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").Value = "'1.97"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'8.36"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").Value = "'143.42"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "2.823.87"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'0.0359"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").Value = "'25.98"
$ws.Range("E51").Value = "  +11.87%  "
